$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.707388877868652
$ws.Range("B1").Value = 2.4711594581604
$ws.Range("C1").Value = 4.661949634552002
$ws.Range("D1").Value = 4.173670768737793
$ws.Range("E1").Value = 1.159024834632874
